$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.936.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.328.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.45%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "182.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "532.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.326.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.29%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.617"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.76%  "

$ws.Range("E12").Value = "  -5.80%  "

$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.837.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.320.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.76%  "

$ws.Range("E17").Value = "  -4.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "64.756.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.968"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.80%  "

$ws.Range("E27").Value = "  -2.34%  "

$ws.Range("E28").Value = "  -4.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "652.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.00%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.396"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0705"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.41%  "

$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.897.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0404"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.89%  "

$ws.Range("E46").Value = "  -3.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.91%  "

$ws.Range("E49").Value = "  +1.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.37%  "
